$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 74667.14
$ws.Range("I64").Value = 3171
$ws.Range("J64").Value = 128289.25
$ws.Range("K64").Value = 3171
$ws.Range("L64").Value = 128289.25
$ws.Range("M64").Value = -2923
$ws.Range("N64").Value = -128785.25
$ws.Range("H67").Value = 74667.14
$ws.Range("I67").Value = 3171
$ws.Range("J67").Value = 128289.25
$ws.Range("K67").Value = 3171
$ws.Range("L67").Value = 128289.25
$ws.Range("M67").Value = -2313
$ws.Range("N67").Value = -130005.25
$ws.Range("H70").Value = 1925.5
$ws.Range("I70").Value = 1751
$ws.Range("J70").Value = 2100
$ws.Range("K70").Value = 5253
$ws.Range("L70").Value = 6300
$ws.Range("M70").Value = -4983
$ws.Range("N70").Value = -6840
$ws.Range("H73").Value = 1925.5
$ws.Range("I73").Value = 1751
$ws.Range("J73").Value = 2100
$ws.Range("K73").Value = 5253
$ws.Range("L73").Value = 6300
$ws.Range("M73").Value = -4317
$ws.Range("N73").Value = -8172
$ws.Range("H74").Value = 3788.8333
$ws.Range("I74").Value = 3794.111
$ws.Range("J74").Value = 3783.5557
$ws.Range("K74").Value = 3794.111
$ws.Range("L74").Value = 3783.5557
$ws.Range("M74").Value = -2858.111
$ws.Range("N74").Value = -5655.5557
$ws.Range("H77").Value = 3788.8333
$ws.Range("I77").Value = 3794.111
$ws.Range("J77").Value = 3783.5557
$ws.Range("K77").Value = 18970.555
$ws.Range("L77").Value = 18917.7785
$ws.Range("M77").Value = -14290.555
$ws.Range("N77").Value = -28277.7785
$ws.Range("H87").Value = 27276
$ws.Range("J87").Value = 27276
$ws.Range("L87").Value = 27276
$ws.Range("N87").Value = -29772
$ws.Range("H90").Value = 27276
$ws.Range("J90").Value = 27276
$ws.Range("L90").Value = 81828
$ws.Range("N90").Value = -94308
$ws.Range("H125").Value = 7516.75
$ws.Range("I125").Value = 26008
$ws.Range("J125").Value = 1353
$ws.Range("K125").Value = 234072
$ws.Range("L125").Value = 12177
$ws.Range("M125").Value = -231612
$ws.Range("N125").Value = -17097
$ws.Range("H137").Value = 2121.9583
$ws.Range("I137").Value = 1236.35
$ws.Range("J137").Value = 6550
$ws.Range("K137").Value = 3709.05
$ws.Range("L137").Value = 19650
$ws.Range("M137").Value = -1159.05
$ws.Range("N137").Value = -24750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9255.846
$ws.Range("I61").Value = 10637.454
$ws.Range("J61").Value = 1657
$ws.Range("K61").Value = 10637.454
$ws.Range("L61").Value = 1657
$ws.Range("M61").Value = -10425.454
$ws.Range("N61").Value = -2081
$ws.Range("H63").Value = 2079.2
$ws.Range("J63").Value = 2396
$ws.Range("L63").Value = 2396
$ws.Range("N63").Value = -3768
$ws.Range("H66").Value = 2079.2
$ws.Range("J66").Value = 2396
$ws.Range("L66").Value = 11980
$ws.Range("N66").Value = -18844
$ws.Range("H98").Value = 26118.334
$ws.Range("J98").Value = 26118.334
$ws.Range("L98").Value = 26118.334
$ws.Range("N98").Value = -32108.334
$ws.Range("H132").Value = 8705.643
$ws.Range("I132").Value = 6107.3
$ws.Range("K132").Value = 18321.9
$ws.Range("M132").Value = -15791.9
$ws.Range("H136").Value = 9255.846
$ws.Range("I136").Value = 10637.454
$ws.Range("J136").Value = 1657
$ws.Range("K136").Value = 31912.362
$ws.Range("L136").Value = 4971
$ws.Range("M136").Value = -29362.362
$ws.Range("N136").Value = -10071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 830.44446
$ws.Range("I64").Value = 720
$ws.Range("J64").Value = 918.8
$ws.Range("K64").Value = 720
$ws.Range("L64").Value = 918.8
$ws.Range("M64").Value = -495
$ws.Range("N64").Value = -1368.8
$ws.Range("H67").Value = 830.44446
$ws.Range("I67").Value = 720
$ws.Range("J67").Value = 918.8
$ws.Range("K67").Value = 720
$ws.Range("L67").Value = 918.8
$ws.Range("M67").Value = 60
$ws.Range("N67").Value = -2478.8
$ws.Range("H86").Value = 3062.3784
$ws.Range("I86").Value = 3008.5557
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3008.5557
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1885.5557
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 3062.3784
$ws.Range("I89").Value = 3008.5557
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 15042.7785
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -9426.7785
$ws.Range("N89").Value = -36232
$ws.Range("H134").Value = 3958.0557
$ws.Range("I134").Value = 4105.758
$ws.Range("J134").Value = 2333.3333
$ws.Range("K134").Value = 12317.274
$ws.Range("L134").Value = 6999.999899999999
$ws.Range("M134").Value = -9782.273999999999
$ws.Range("N134").Value = -12069.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4005.805
$ws.Range("I31").Value = 5035.6665
$ws.Range("J31").Value = 2551.8823
$ws.Range("K31").Value = 5035.6665
$ws.Range("L31").Value = 2551.8823
$ws.Range("M31").Value = -4740.6665
$ws.Range("N31").Value = -3141.8823
$ws.Range("H34").Value = 4005.805
$ws.Range("I34").Value = 5035.6665
$ws.Range("J34").Value = 2551.8823
$ws.Range("K34").Value = 5035.6665
$ws.Range("L34").Value = 2551.8823
$ws.Range("M34").Value = -4833.6665
$ws.Range("N34").Value = -2955.8823
$ws.Range("H62").Value = 5851.5
$ws.Range("I62").Value = 2400
$ws.Range("J62").Value = 7002
$ws.Range("K62").Value = 2400
$ws.Range("L62").Value = 7002
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -8250
$ws.Range("H65").Value = 5851.5
$ws.Range("I65").Value = 2400
$ws.Range("J65").Value = 7002
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 35010
$ws.Range("M65").Value = -8880
$ws.Range("N65").Value = -41250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1728.9166
$ws.Range("I132").Value = 987.5
$ws.Range("J132").Value = 2099.625
$ws.Range("K132").Value = 8887.5
$ws.Range("L132").Value = 18896.625
$ws.Range("M132").Value = -6357.5
$ws.Range("N132").Value = -23956.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5055.575
$ws.Range("I70").Value = 5001
$ws.Range("J70").Value = 5182.9165
$ws.Range("K70").Value = 5001
$ws.Range("L70").Value = 5182.9165
$ws.Range("M70").Value = -4731
$ws.Range("N70").Value = -5722.9165
$ws.Range("H73").Value = 5055.575
$ws.Range("I73").Value = 5001
$ws.Range("J73").Value = 5182.9165
$ws.Range("K73").Value = 5001
$ws.Range("L73").Value = 5182.9165
$ws.Range("M73").Value = -4065
$ws.Range("N73").Value = -7054.9165
$ws.Range("H80").Value = 3053.077
$ws.Range("I80").Value = 2956.25
$ws.Range("J80").Value = 3208
$ws.Range("K80").Value = 2956.25
$ws.Range("L80").Value = 3208
$ws.Range("M80").Value = -1958.25
$ws.Range("N80").Value = -5204
$ws.Range("H83").Value = 3053.077
$ws.Range("I83").Value = 2956.25
$ws.Range("J83").Value = 3208
$ws.Range("K83").Value = 14781.25
$ws.Range("L83").Value = 16040
$ws.Range("M83").Value = -9789.25
$ws.Range("N83").Value = -26024

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 13500
$ws.Range("J106").Value = 13500
$ws.Range("L106").Value = 13500
$ws.Range("N106").Value = -16024

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 14320
$ws.Range("J64").Value = 14320
$ws.Range("L64").Value = 14320
$ws.Range("N64").Value = -14816
$ws.Range("H67").Value = 14320
$ws.Range("J67").Value = 14320
$ws.Range("L67").Value = 14320
$ws.Range("N67").Value = -16036
